# Edit: insert a new data row at row 422 (pushing all subsequent rows down by
# one), and populate the new row with a fresh record. This matches the
# target diff, where the whole data block from row 422 through the former
# last row 509 shifts down by one row, and the last original row (509)
# becomes new row 510.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 422; Excel shifts rows 422:509 down to 423:510,
# carrying their formatting (including the date style on column D) with them.
$ws.Rows("422").Insert()

# Populate the newly inserted row 422 with the new record's values.
$ws.Range("A422").Value2 = 6
$ws.Range("B422").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C422").Value2 = "Metropolitana"
$ws.Range("D422").Value2 = 44637
$ws.Range("E422").Value2 = 13
$ws.Range("F422").Value2 = 100112012
$ws.Range("G422").Value2 = "Espinaca"
$ws.Range("H422").Value2 = "Sin especificar"
$ws.Range("I422").Value2 = "Primera"
$ws.Range("J422").Value2 = 530
$ws.Range("K422").Value2 = 6000
$ws.Range("L422").Value2 = 7000
$ws.Range("M422").Value2 = 6453
$ws.Range("N422").Value2 = '$/cuna 10 kilos'
$ws.Range("O422").Value2 = "Región Metropolitana"
$ws.Range("P422").Value2 = 645
$ws.Range("Q422").Value2 = 10
$ws.Range("R422").Value2 = "Hortaliza"
